$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, column index, new value
$updates = @(
    @(2, 7, 0.05931326035507824),
    @(2, 8, -7.758583532447306),
    @(2, 9, -11.68258649709265),
    @(3, 7, 0.06787059382195859),
    @(3, 8, 20.75105541663314),
    @(4, 7, -0.01721485920497215),
    @(4, 8, 15.24620223976365),
    @(5, 7, -0.01250155079535128),
    @(5, 8, -9.362162939859571),
    @(6, 7, -0.01102466719485852),
    @(6, 8, 1.508818547482636),
    @(7, 7, -0.0181440260750904),
    @(7, 8, -220.3543745543151),
    @(8, 7, 0.003502734656408314),
    @(8, 8, 161.5685410690318),
    @(9, 7, 0.002892286301756805),
    @(9, 8, 152.6526361526776),
    @(10, 7, -0.07058793924194762),
    @(10, 8, -12.00567110638592),
    @(11, 7, -0.0717485056820435),
    @(11, 8, -11.97049213444053),
    @(12, 7, -0.4077367820469791),
    @(12, 8, -3.328248355214717),
    @(13, 7, -0.4003244468899067),
    @(13, 8, -2.092538887822796),
    @(14, 7, -0.01476982522068948),
    @(14, 8, -80.24591462700954),
    @(15, 7, -0.02049596885707215),
    @(15, 8, 54.77368071823023),
    @(16, 7, 0.1357106653806431),
    @(16, 8, -0.7597708373962435),
    @(17, 7, 0.1419156419427633),
    @(17, 8, 1.753059606498887),
    @(18, 7, 0.1183490751769825),
    @(18, 8, 0.5786332321600008),
    @(19, 7, 0.1276814192229102),
    @(19, 8, -0.7674211773128655),
    @(20, 7, 0.08466092869781687),
    @(20, 8, -4.590883935027265),
    @(21, 7, 0.08701871041812406),
    @(21, 8, -0.08303364102259785),
    @(22, 7, -0.09250295906055113),
    @(22, 8, 1.048542055660304),
    @(23, 7, -0.1044973442376267),
    @(23, 8, -3.008847718447412),
    @(24, 7, 0.1602526855368297),
    @(24, 8, -0.5204360456052857),
    @(25, 7, 0.1729861740892776),
    @(25, 8, 1.401232081710936),
    @(26, 7, 0.08703935028236293),
    @(26, 8, -3.987086673039802),
    @(27, 7, 0.08452875962460425),
    @(27, 8, -1.667005053999603),
    @(28, 7, -0.1401135276862405),
    @(28, 8, -1.796776725288932),
    @(29, 7, -0.1418216398868193),
    @(29, 8, -1.446264527913983),
    @(30, 7, 0.04980383740963305),
    @(30, 8, -4.251514958609911),
    @(31, 7, 0.03986722474128181),
    @(31, 8, -9.019242618376236),
    @(32, 7, 0.110250098306307),
    @(32, 8, 1.422364966319449),
    @(33, 7, 0.1125325945037176),
    @(33, 8, -9.317522868537912),
    @(34, 7, -0.01720642034022659),
    @(34, 8, -10.18177691115705),
    @(35, 7, -0.02032978017675743),
    @(35, 8, -21.46400111070656),
    @(36, 7, 0.03102358377163858),
    @(36, 8, -15.62090973897721),
    @(37, 7, 0.03905474384722487),
    @(37, 8, 9.441037305167402),
    @(38, 7, 0.09524104742977489),
    @(38, 8, -5.045931657436529),
    @(39, 7, 0.1029960564212245),
    @(39, 8, 5.745421552887039),
    @(40, 7, 0.03181358538296872),
    @(40, 8, -5.562976482689675),
    @(41, 7, 0.03135829096378918),
    @(41, 8, -2.670602482397569),
    @(42, 7, 0.1195023676515497),
    @(42, 8, -1.160463039777635),
    @(43, 7, 0.1196015280213674),
    @(43, 8, -6.404421917759387),
    @(44, 7, 0.03709734698112794),
    @(44, 8, -6.477211797802377),
    @(45, 7, 0.03005657039346328),
    @(45, 8, -3.548456498734716),
    @(46, 7, 0.05970021648889693),
    @(46, 8, 5.444504780484835),
    @(47, 7, 0.06451399839426483),
    @(47, 8, 9.960445098339745),
    @(48, 7, 0.04508408040356845),
    @(48, 8, -8.472328758822391),
    @(49, 7, 0.04992211562900956),
    @(49, 8, 9.555812187511888),
    @(50, 7, 0.02768077958099759),
    @(50, 8, 4.509076548496136),
    @(51, 7, 0.02779305827376145),
    @(51, 8, -0.7952999442275415),
    @(52, 7, -0.08628913714614755),
    @(52, 8, 0.7299555758046753),
    @(53, 7, -0.08183443460445265),
    @(53, 8, -2.017338937103921),
    @(54, 7, 0.04280035091415067),
    @(54, 8, -14.44498784602836),
    @(55, 7, 0.04763785654799152),
    @(55, 8, -15.36322556036827),
    @(56, 7, 0.04197547155325531),
    @(56, 8, -15.08423546382736),
    @(57, 7, 0.04572643212521518),
    @(57, 8, 20.38347404259697),
    @(58, 7, 0.05952686678328963),
    @(58, 8, 3.320047951676589),
    @(59, 7, 0.06487831750109789),
    @(59, 8, 13.79503894456258),
    @(60, 7, 0.0261358622762431),
    @(60, 8, -4.822058504832656),
    @(61, 7, 0.0280465791681736),
    @(61, 8, 5.05192402506522),
    @(62, 7, 0.05839014763650523),
    @(62, 8, -6.496819078686856),
    @(63, 7, 0.06799718271996416),
    @(63, 8, 6.440734877502613),
    @(64, 7, 0.02745762061315563),
    @(64, 8, -1.01672780668489),
    @(65, 7, 0.03573320595517088),
    @(65, 8, 0.8636121236346772),
    @(66, 7, 0.08152095606484824),
    @(66, 8, 4.940290602668949),
    @(67, 7, 0.07517468216269517),
    @(67, 8, -4.686334403711003),
    @(68, 7, -0.02984344717803936),
    @(68, 8, -37.23842943917121),
    @(69, 7, -0.01054679920028078),
    @(69, 8, 44.90033341459478),
    @(70, 7, 0.07493595258601712),
    @(70, 8, 4.098151554425201),
    @(71, 7, 0.06735632600752715),
    @(71, 8, -15.19082621697576),
    @(72, 7, -0.1499479677930243),
    @(72, 8, 2.403982637671619),
    @(73, 7, -0.1561318177261071),
    @(73, 8, -1.996571493005017),
    @(74, 7, 0.1450107745304179),
    @(74, 8, -3.595564635397051),
    @(75, 7, 0.1617382953170958),
    @(75, 8, 7.504932914373197),
    @(76, 7, -0.01067521366602649),
    @(76, 8, -929.9670582768416),
    @(77, 7, -0.003360713474036557),
    @(77, 8, -52.2139619980198),
    @(78, 7, 0.0942621982319341),
    @(78, 8, 4.784383639568476),
    @(79, 7, 0.09253109635554063),
    @(79, 8, -4.510017358421345),
    @(80, 7, -0.2150662454157058),
    @(80, 8, 0.6316597404747222),
    @(81, 7, -0.2048417367402791),
    @(81, 8, 3.881562475961414),
    @(82, 7, 0.1686638702000806),
    @(82, 8, 0.6241597360466682),
    @(83, 7, 0.183976288312983),
    @(83, 8, 4.514801078239142),
    @(84, 7, 0.111158268028414),
    @(84, 8, 4.738137537994417),
    @(85, 7, 0.1171861838117136),
    @(85, 8, 12.07681808993687)
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}
